$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.683.47"
$ws.Range("E2").Value = "  +0.88%  "

$ws.Range("D3").Value = "1.887.75"
$ws.Range("E3").Value = "  +1.00%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'247.80"
$ws.Range("E5").Value = "  +0.53%  "

$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").Value = "'0.4730"
$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("E8").Value = "  +0.62%  "

$ws.Range("D9").Value = "'0.06527"
$ws.Range("E9").Value = "  +0.47%  "

$ws.Range("D10").Value = "'22.07"
$ws.Range("E10").Value = "  +0.37%  "

$ws.Range("D11").Value = "'0.07809"
$ws.Range("E11").Value = "  +1.19%  "

$ws.Range("D12").Value = "1.889.42"
$ws.Range("E12").Value = "  +1.08%  "

$ws.Range("D13").Value = "'96.76"
$ws.Range("E13").Value = "  -0.74%  "

$ws.Range("D14").Value = "'0.7369"
$ws.Range("E14").Value = "  -0.18%  "

$ws.Range("E15").Value = "  +2.65%  "

$ws.Range("D16").Value = "'284.30"
$ws.Range("E16").Value = "  +3.74%  "

$ws.Range("D17").Value = "30.669.75"
$ws.Range("E17").Value = "  +0.88%  "

$ws.Range("D18").Value = "'13.26"
$ws.Range("E18").Value = "  -0.71%  "

$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").Value = "'0.000007528"

$ws.Range("D21").Value = "2.138.15"

$ws.Range("D22").Value = "'5.316"
$ws.Range("E22").Value = "  +1.78%  "

$ws.Range("D23").Value = "'0.9993"
$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("D24").Value = "'6.245"
$ws.Range("E24").Value = "  +1.33%  "

$ws.Range("D25").Value = "'9.225"
$ws.Range("E25").Value = "  -0.56%  "

$ws.Range("D26").Value = "'164.87"
$ws.Range("E26").Value = "  +0.40%  "

$ws.Range("D27").Value = "'18.95"
$ws.Range("E27").Value = "  +0.67%  "

$ws.Range("D28").Value = "'1.917"
$ws.Range("E28").Value = "  -0.64%  "

$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").Value = "'0.09744"
$ws.Range("E29").Value = "  -2.53%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'1.335"
$ws.Range("E30").Value = "  -2.19%  "

$ws.Range("D31").Value = "'1.488"
$ws.Range("E31").Value = "  -1.21%  "

$ws.Range("D32").Value = "'4.298"
$ws.Range("E32").Value = "  -0.14%  "

$ws.Range("D33").Value = "'4.192"
$ws.Range("E33").Value = "  +1.02%  "

$ws.Range("D34").Value = "'0.04864"
$ws.Range("E34").Value = "  +0.62%  "

$ws.Range("D35").Value = "'1.125"
$ws.Range("E35").Value = "  +0.57%  "

$ws.Range("D36").Value = "'0.6963"
$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("E37").Value = "  +0.46%  "

$ws.Range("D38").Value = "'0.01894"
$ws.Range("E38").Value = "  +1.94%  "

$ws.Range("D39").Value = "'2.807"
$ws.Range("E39").Value = "  +2.27%  "

$ws.Range("D40").Value = "'6.351"
$ws.Range("E40").Value = "  +0.82%  "

$ws.Range("D41").Value = "'76.17"
$ws.Range("E41").Value = "  +4.56%  "

$ws.Range("D42").Value = "'1.996"
$ws.Range("E42").Value = "  +1.37%  "

$ws.Range("D43").Value = "'0.4268"
$ws.Range("E43").Value = "  +1.87%  "

$ws.Range("E44").Value = "  +0.09%  "

$ws.Range("D45").Value = "'0.8363"
$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("D46").Value = "'101.62"
$ws.Range("E46").Value = "  -0.46%  "

$ws.Range("D47").Value = "'9.483"
$ws.Range("E47").Value = "  +2.73%  "

$ws.Range("D48").Value = "'35.59"
$ws.Range("E48").Value = "  +0.83%  "

$ws.Range("D49").Value = "'7.031"
$ws.Range("E49").Value = "  +0.28%  "

$ws.Range("D50").Value = "'915.19"

$ws.Range("E51").Value = "  +2.12%  "
